$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 171 - this shifts the existing rows
# 171..235 down to 172..236 (keeping all their data untouched), and
# the sheet's used-range grows from A1:R235 to A1:R236.
$ws.Rows("171:171").Insert()

# Populate the newly inserted row 171 with the new record.
$ws.Range("A171").Value2 = 10
$ws.Range("B171").Value2 = "Vega Modelo de Temuco"
$ws.Range("C171").Value2 = "La Araucanía"
$ws.Range("D171").Value2 = 44809
$ws.Range("E171").Value2 = 9
$ws.Range("F171").Value2 = 100112005
$ws.Range("G171").Value2 = "Puerro"
$ws.Range("H171").Value2 = "Azul de Maquehue"
$ws.Range("I171").Value2 = "Primera"
$ws.Range("J171").Value2 = 30
$ws.Range("K171").Value2 = 20000
$ws.Range("L171").Value2 = 20000
$ws.Range("M171").Value2 = 20000
$ws.Range("N171").Value2 = "$/docena de paquetes"
$ws.Range("O171").Value2 = "Provincia de Cautín"
$ws.Range("P171").Value2 = 1667
$ws.Range("Q171").Value2 = 12
$ws.Range("R171").Value2 = "Hortaliza"
